# edit.ps1 -- applies the E_EV05_Perfil_Emprendedor.docx revision
# described by the supplied unified diff, via Word COM interop.
#
# Summary of changes applied:
#  1. Body text: "...en tu país. En el caso de Colombia, hay leyes..."
#     -> "...en Colombia, hay leyes..."
#  2. Body text: "Después, debes conocer..." is retyped as three runs
#     "Después, " + "se debe " + "conocer..." (identical formatting)
#     and "¿O es una oportunidad" -> "¿ una oportunidad".
#  3. Hyperlink run "_ Perfil_Emprendedor" is split into two runs
#     "_ " and "Perfil_Emprendedor" (identical formatting/style).
#
# (Cosmetic, save-time-only artifacts that Word's UI/engine produce on
# their own -- <w:lastRenderedPageBreak/>, <w:proofErr/> spell-check
# markers, <w:noProof/> on the inserted picture run, and the internal
# renumbering of the customXml parts -- have no corresponding property
# in the Word object model exposed by this host, so they are not
# reproduced here.)

# Word normalises/merges adjacent runs that end up with identical
# formatting whenever their *text* is touched. To keep two
# same-formatted runs apart (matching runs Word itself would leave
# behind from separate edits) we flip a character property off then
# back to its original value; this is a no-op for the rendered
# formatting but forces the run boundary to persist.
function Force-RunBoundary($rng) {
    $orig = $rng.Font.Bold
    if ($orig -eq 1 -or $orig -eq -1) {
        $rng.Font.Bold = 0
        $rng.Font.Bold = $orig
    } else {
        $rng.Font.Bold = 1
        $rng.Font.Bold = $orig
    }
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "en tu país. En el caso de Colombia" -> "en Colombia"
# ---------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "en tu país. En el caso de Colombia, hay leyes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "en Colombia, hay leyes", 2)
Write-Output ("Edit1 (reglas del juego): " + $found1)

# ---------------------------------------------------------------------
# Edit 2: "Después, debes conocer" -> "Después, se debe conocer"
#         and "¿O es una oportunidad" -> "¿ una oportunidad"
# ---------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "Después, debes conocer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Después, se debe conocer", 2)
Write-Output ("Edit2a (se debe): " + $found2)

$found3 = $d.Content.Find.Execute(
    "¿O es una oportunidad",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "¿ una oportunidad", 2)
Write-Output ("Edit2b (oportunidad): " + $found3)

# Re-locate the paragraph (text has changed, but the paragraph itself
# has not been added/removed, so the index is stable) and split it
# into three runs:
#   "Después, " | "se debe " | "conocer qué tipo de emprendimiento..."
$p = $d.Paragraphs(13)
$pStart = $p.Range.Start

$seg1 = $d.Range($pStart, $pStart + 9)          # "Después, "
Write-Output ("seg1: [" + $seg1.Text + "]")
Force-RunBoundary($seg1)

$seg2 = $d.Range($pStart + 9, $pStart + 9 + 8)  # "se debe "
Write-Output ("seg2: [" + $seg2.Text + "]")
Force-RunBoundary($seg2)

# ---------------------------------------------------------------------
# Edit 3: split the hyperlink run "_ Perfil_Emprendedor" into
#         "_ " and "Perfil_Emprendedor"
# ---------------------------------------------------------------------
$hr = $d.Content
$foundH = $hr.Find.Execute(
    "_ Perfil_Emprendedor", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
Write-Output ("Edit3 locate: " + $foundH)

$hStart = $hr.Start
$hseg1 = $d.Range($hStart, $hStart + 2)   # "_ "
Write-Output ("hseg1: [" + $hseg1.Text + "]")
Force-RunBoundary($hseg1)

Write-Output "Done."
